# Prepare project for Render deployment
# Populate the "Enquiries" sheet with the first captured enquiry-form
# submission: a header row followed by the submitted record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-existing boilerplate number format (carried over from the export
# template) applied across the used range before the data is entered.
$ws.Range("A1:D2").NumberFormat = """上午/下午 ""hh""時""mm""分""ss""秒 """

# Header row
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "Phone"
$ws.Range("D1").Value = "Message"

# Enquiry data row
$ws.Range("A2").Value = "Harsh Gahlot"
$ws.Range("B2").Value = "gahlotharsh19@gmail.com"

# Phone number must stay text so the leading zero survives
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "09760978001"

$ws.Range("D2").Value = "kjkjadakdjakfjadkjfas;jf;adsjkjal;dsn;lkdjasdk;as"

# Drop the now-unused trailing blank rows that used to pad the sheet
$null = $ws.Rows("3:6").Delete()
